$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The test app changed (now https://opencart.abstracta.us/index.php), so the
# stored login credentials in row 2 (A2/B2) need to be updated accordingly.
$ws.Range("A2").Value = "Test7788@gmail.com"
$ws.Range("B2").Value = "Test7788"

# Re-fit the columns that actually hold data to their (now different) text,
# just like Excel does via Home > Format > AutoFit Column Width.
foreach ($col in @("A", "B", "C", "E", "G", "H")) {
    $ws.Columns($col).AutoFit()
}

# The workbook was left with the selection on B3.
$ws.Range("B3").Select()
